# The edit swaps the data of row 58 and row 60 (the "Id" record 112145544 /
# Knärot entry and the "Id" record 112145535 / Talltita entry trade places
# in the sheet), while a handful of columns that already held identical
# values in both rows (C, K, L, N, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX,
# AY) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: reading via the parameterized `.Value` getter doesn't resolve in
# this host, so use `.Value2` (plain property) for all reads; `.Value =`
# remains fine for writes.

# --- capture the "before" values of row 58 ---
$A58 = $ws.Range("A58").Value2
$B58 = $ws.Range("B58").Value2
$D58 = $ws.Range("D58").Value2
$E58 = $ws.Range("E58").Value2
$F58 = $ws.Range("F58").Value2
$G58 = $ws.Range("G58").Value2
$H58 = $ws.Range("H58").Value2
$I58 = $ws.Range("I58").Value2
$J58 = $ws.Range("J58").Value2
$M58 = $ws.Range("M58").Value2
$P58 = $ws.Range("P58").Value2
$Q58 = $ws.Range("Q58").Value2
$R58 = $ws.Range("R58").Value2
$S58 = $ws.Range("S58").Value2
$Z58 = $ws.Range("Z58").Value2
$AB58 = $ws.Range("AB58").Value2

# --- capture the "before" values of row 60 ---
$A60 = $ws.Range("A60").Value2
$B60 = $ws.Range("B60").Value2
$D60 = $ws.Range("D60").Value2
$E60 = $ws.Range("E60").Value2
$F60 = $ws.Range("F60").Value2
$G60 = $ws.Range("G60").Value2
$H60 = $ws.Range("H60").Value2
$I60 = $ws.Range("I60").Value2
$M60 = $ws.Range("M60").Value2
$P60 = $ws.Range("P60").Value2
$Q60 = $ws.Range("Q60").Value2
$R60 = $ws.Range("R60").Value2
$S60 = $ws.Range("S60").Value2
$Z60 = $ws.Range("Z60").Value2
$AB60 = $ws.Range("AB60").Value2

# Column I ("Antal") is stored as *text* in this sheet even though the
# values look numeric (e.g. "17"), so force text formatting before the
# write to stop the host from auto-coercing the digit string to a number,
# then restore the default style so no stray number format sticks around.

# --- write row 60's old values into row 58 ---
$ws.Range("A58").Value = $A60
$ws.Range("B58").Value = $B60
$ws.Range("D58").Value = $D60
$ws.Range("E58").Value = $E60
$ws.Range("F58").Value = $F60
$ws.Range("G58").Value = $G60
$ws.Range("H58").Value = $H60
$ws.Range("I58").NumberFormat = "@"
$ws.Range("I58").Value = $I60
$ws.Range("I58").Style = "Normal"
$ws.Range("J58").ClearContents()
$ws.Range("M58").Value = $M60
$ws.Range("P58").Value = $P60
$ws.Range("Q58").Value = $Q60
$ws.Range("R58").Value = $R60
$ws.Range("S58").Value = $S60
$ws.Range("Z58").Value = $Z60
$ws.Range("AB58").Value = $AB60

# --- write row 58's old values into row 60 ---
$ws.Range("A60").Value = $A58
$ws.Range("B60").Value = $B58
$ws.Range("D60").Value = $D58
$ws.Range("E60").Value = $E58
$ws.Range("F60").Value = $F58
$ws.Range("G60").Value = $G58
$ws.Range("H60").Value = $H58
$ws.Range("I60").NumberFormat = "@"
$ws.Range("I60").Value = $I58
$ws.Range("I60").Style = "Normal"
$ws.Range("J60").Value = $J58
$ws.Range("M60").ClearContents()
$ws.Range("P60").Value = $P58
$ws.Range("Q60").Value = $Q58
$ws.Range("R60").Value = $R58
$ws.Range("S60").Value = $S58
$ws.Range("Z60").Value = $Z58
$ws.Range("AB60").Value = $AB58
